$wb = $excel.ActiveWorkbook

# Add new "Tasks" worksheet and move it to be the last tab.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Tasks"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch a stable reference to the sheet by name for all further edits.
$ts = $wb.Worksheets.Item("Tasks")

$ts.Range("A1").Value = "Title"
$ts.Range("B1").Value = "Completion"
$ts.Range("A2").Value = "BVC "
$ts.Range("B2").Value = 25
$ts.Range("A3").Value = "ATCO Lab"
$ts.Range("B3").Value = 65
$ts.Range("A4").Value = "General Dynamics"
$ts.Range("B4").Value = 10

$ts.Columns.Item(1).ColumnWidth = 17.5
$ts.Columns.Item(2).ColumnWidth = 16

# Make Tasks the active/selected sheet with D3 selected (matches target tab state).
$ts.Activate()
$ts.Range("D3").Select()
